$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the start/end time entries for the row dated 2020-10-15 (row 19)
$ws.Range("B19").Value = 0.47355324074074073
$ws.Range("C19").Value = 0.6651273148148148

# Add the new changelog note text to column F for that row, matching the
# formatting already used by the other note cells in that column.
$ws.Range("F18").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newNote = "* Added NPC node`n* Added node connection validation`n* Added styling and custom node layouts (UIElements learning)`n* Fixed issue where renaming a property would not rename nodes that refer to that property`n* Started working on copy&paste functionality inside the tool"
$ws.Range("F19").Value = $newNote

# Keep the row height as the other (non-autofit) rows.
$ws.Rows.Item(19).RowHeight = 15
$ws.Rows.Item(19).EntireRow.AutoFit() | Out-Null
$ws.Rows.Item(19).RowHeight = 15

# Update the active selection to G19
$ws.Range("G19").Select()
